$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'69.739.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.02%  "
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'3.491.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.99%  "
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.23%  "
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'606.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +3.27%  "
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'192.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +2.43%  "
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'0.625"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.75%  "
$ws.Range("E7").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'0.213"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.91%  "
$ws.Range("E9").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'53.31"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.62%  "
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.0000306"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.55%  "
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'9.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.79%  "
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'4.051.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.90%  "
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").Value = "'608.07"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +4.95%  "
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'69.785.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.00%  "
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'12.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.01%  "
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'18.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.75%  "
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'3.495.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.43%  "
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("E20").Value = "'  -0.13%  "
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("E21").Value = "'  -0.79%  "
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'17.80"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.29%  "
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'105.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +11.74%  "
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("E24").Value = "'  -0.01%  "
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'5.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +3.96%  "
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = "'  +4.35%  "
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("E27").Value = "'  -0.41%  "
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("D28").Value = "'9.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +5.31%  "
$ws.Range("E28").Style = "Normal"

# Row 29
$ws.Range("D29").Value = "'34.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +5.47%  "
$ws.Range("E29").Style = "Normal"

# Row 30
$ws.Range("D30").Value = "'7.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +1.05%  "
$ws.Range("E30").Style = "Normal"

# Row 31
$ws.Range("D31").Value = "'4.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +14.34%  "
$ws.Range("E31").Style = "Normal"

# Row 32
$ws.Range("D32").Value = "'12.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +3.71%  "
$ws.Range("E32").Style = "Normal"

# Row 33
$ws.Range("D33").Value = "'0.115"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.03%  "
$ws.Range("E33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'64.18"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.46%  "
$ws.Range("E34").Style = "Normal"

# Row 35
$ws.Range("D35").Value = "'3.716.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -1.53%  "
$ws.Range("E35").Style = "Normal"

# Row 36
$ws.Range("D36").Value = "'0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.18%  "
$ws.Range("E36").Style = "Normal"

# Row 37
$ws.Range("B37").Value = "'Bittensor"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'519.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.92%  "
$ws.Range("E37").Style = "Normal"

# Row 38
$ws.Range("B38").Value = "'Fetch.AI"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'3.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -5.59%  "
$ws.Range("E38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.0₃0793"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.72%  "
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'3.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.06%  "
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("E41").Value = "'  -4.14%  "
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'36.52"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -4.03%  "
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("E43").Value = "'  -1.07%  "
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("E44").Value = "'  +1.27%  "
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("E45").Value = "'  -3.28%  "
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("E46").Value = "'  +1.65%  "
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("E47").Value = "'  -4.00%  "
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("E48").Value = "'  +0.46%  "
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("E49").Value = "'  -5.27%  "
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'133.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.74%  "
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'1.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +10.42%  "
$ws.Range("E51").Style = "Normal"

